$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "Armor" is added after new strings for row 6 so that
# shared-string order comes out as: ... Haunty Squire(8), Protector(9), Armor(10)
$ws.Range("C6").Value = "Minion"
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = "Protector"
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 30
$ws.Range("H3").Value = "Armor"
$ws.Range("H6").Value = 10

# Update the active selection to match the post-edit state (B6 selected)
$ws.Range("B6").Select() | Out-Null
